$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2 through 43 changes from 45831 to 45832
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 3).Value = 45832
}
